$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sheet1")

$ws.Range("D2").Value = '26.077.32'
$ws.Range("E2").Value = '  +0.37%  '
$ws.Range("D3").Value = '1.639.46'
$ws.Range("E3").Value = '  -0.10%  '
$ws.Range("E4").Value = '  +0.46%  '
$ws.Range("D5").Value = "'214.42"
$ws.Range("D5").Style = "Normal"
$ws.Range("E5").Value = '  -0.69%  '
$ws.Range("E6").Value = '  -0.25%  '
$ws.Range("E7").Value = '  +0.47%  '
$ws.Range("E8").Value = '  -2.45%  '
$ws.Range("E9").Value = '  -2.40%  '
$ws.Range("D10").Value = "'18.60"
$ws.Range("D10").Style = "Normal"
$ws.Range("E10").Value = '  -5.27%  '
$ws.Range("E11").Value = '  -0.19%  '
$ws.Range("D12").Value = '1.750.87'
$ws.Range("E12").Value = '  +6.22%  '
$ws.Range("D13").Value = "'4.21"
$ws.Range("D13").Style = "Normal"
$ws.Range("E13").Value = '  -1.67%  '
$ws.Range("D14").Value = "'0.530"
$ws.Range("D14").Style = "Normal"
$ws.Range("E14").Value = '  -2.76%  '
$ws.Range("D15").Value = "'62.32"
$ws.Range("D15").Style = "Normal"
$ws.Range("E15").Value = '  -1.17%  '
$ws.Range("E16").Value = '  -2.05%  '
$ws.Range("D17").Value = '26.081.58'
$ws.Range("E17").Value = '  +0.54%  '
$ws.Range("E18").Value = '  +0.51%  '
$ws.Range("D19").Value = "'190.38"
$ws.Range("D19").Style = "Normal"
$ws.Range("E19").Value = '  -1.38%  '
$ws.Range("E20").Value = '  -2.11%  '
$ws.Range("D21").Value = "'9.57"
$ws.Range("D21").Style = "Normal"
$ws.Range("E21").Value = '  -3.61%  '
$ws.Range("E22").Value = '  -2.60%  '
$ws.Range("D23").Value = "'144.25"
$ws.Range("D23").Style = "Normal"
$ws.Range("E23").Value = '  +0.38%  '
$ws.Range("E24").Value = '  -0.47%  '
$ws.Range("E25").Value = '  +0.42%  '
$ws.Range("E26").Value = '  -1.60%  '
$ws.Range("D27").Value = "'6.77"
$ws.Range("D27").Style = "Normal"
$ws.Range("E27").Value = '  -1.49%  '
$ws.Range("D28").Value = "'15.23"
$ws.Range("D28").Style = "Normal"
$ws.Range("E28").Value = '  -2.42%  '
$ws.Range("E29").Value = '  -0.66%  '
$ws.Range("D30").Value = "'0.0485"
$ws.Range("D30").Style = "Normal"
$ws.Range("E30").Value = '  -3.58%  '
$ws.Range("E31").Value = '  -2.32%  '
$ws.Range("D32").Value = "'3.18"
$ws.Range("D32").Style = "Normal"
$ws.Range("E32").Value = '  -3.79%  '
$ws.Range("D33").Value = "'2.45"
$ws.Range("D33").Style = "Normal"
$ws.Range("E33").Value = '  -0.32%  '
$ws.Range("E34").Value = '  -1.85%  '
$ws.Range("D35").Value = "'0.879"
$ws.Range("D35").Style = "Normal"
$ws.Range("E35").Value = '  -2.48%  '
$ws.Range("D36").Value = '1.122.58'
$ws.Range("E36").Value = '  -1.16%  '
$ws.Range("E37").Value = '  -0.15%  '
$ws.Range("D38").Value = "'0.523"
$ws.Range("D38").Style = "Normal"
$ws.Range("E38").Value = '  -3.96%  '
$ws.Range("E39").Value = '  -1.57%  '
$ws.Range("D40").Value = "'98.80"
$ws.Range("D40").Style = "Normal"
$ws.Range("E40").Value = '  -0.64%  '
$ws.Range("D41").Value = "'0.787"
$ws.Range("D41").Style = "Normal"
$ws.Range("E41").Value = '  -1.51%  '
$ws.Range("E42").Value = '  -3.67%  '
$ws.Range("E43").Value = '  +0.18%  '
$ws.Range("D44").Value = "'55.23"
$ws.Range("D44").Style = "Normal"
$ws.Range("E44").Value = '  -2.65%  '
$ws.Range("E45").Value = '  -1.84%  '
$ws.Range("E46").Value = '  -0.38%  '
$ws.Range("E47").Value = '  -0.02%  '
$ws.Range("D48").Value = "'7.57"
$ws.Range("D48").Style = "Normal"
$ws.Range("E48").Value = '  -1.53%  '
$ws.Range("E49").Value = '  +0.25%  '
$ws.Range("D50").Value = "'0.0928"
$ws.Range("D50").Style = "Normal"
$ws.Range("E50").Value = '  -3.66%  '
$ws.Range("E51").Value = '  -1.16%  '
